# Allow WaterOilGas construction for two-phase problems.
# The "kroend" parameter is split into two distinct parameters:
#   - krogend (kept in the original column, renamed)
#   - krowend (new column inserted right after krogend)
# so that the two-phase (oil relative permeability end-point) values for the
# gas-oil and water-oil systems can be specified independently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column right after the existing "kroend" column (column V),
# shifting swirr/a/b/poro_ref/perm_ref/drho one column to the right.
$ws.Columns("V").Insert() | Out-Null

# Rename the old "kroend" header to "krogend" and label the newly
# inserted column "krowend".
$ws.Range("U1").Value = "krogend"
$ws.Range("V1").Value = "krowend"

# Populate the new krowend column with its value for every data row
# (rows 2 through 10).
for ($r = 2; $r -le 10; $r++) {
    $ws.Range("V$r").Value = 1
}

# Update the active selection to match the edited workbook state.
$ws.Range("T7").Select() | Out-Null
